{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Find the start of the section to remove: the \"Kock\u00e1zat kezel\u00e9s\u00e9nek m\u00f3dszerei\"\n// heading paragraph. Everything from there to the end of the document gets\n// deleted (this whole trailing section, added later, is being reverted).\nconst items = paragraphs.items;\nlet startIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.trim() === \"Kock\u00e1zat kezel\u00e9s\u00e9nek m\u00f3dszerei\") {\n    startIndex = i;\n    break;\n  }\n}\n\nif (startIndex !== -1) {\n  // Delete paragraphs from the end backwards so earlier indices stay valid.\n  for (let i = items.length - 1; i >= startIndex; i--) {\n    items[i].delete();\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the start of the \"Kock\u00e1zat kezel\u00e9s\u00e9nek m\u00f3dszerei\" section (the\n# heading that begins the whole block being removed) and delete everything\n# from there through the end of the document body.\n$range = $d.Content\n$range.Find.ClearFormatting()\n$found = $range.Find.Execute(\"Kock\u00e1zat kezel\u00e9s\u00e9nek m\u00f3dszerei\")\n\nif ($found) {\n    $startPos = $range.Paragraphs.Item(1).Range.Start\n    $delRange = $d.Range($startPos, $d.Content.End)\n    $delRange.Delete()\n}\n"}
